# Update player data: swap out several players for new ones and adjust
# their associated stats (rating, position, country, league, club, cost,
# chemistry) to add additional objective types.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2: Bixby -> Shaw
$ws.Range("A2").Value = "Shaw"
$ws.Range("B2").Value = 83
$ws.Range("E2").Value = "CB"
$ws.Range("G2").Value = "England"
$ws.Range("H2").Value = "Premier League"
$ws.Range("I2").Value = "Manchester Utd"
$ws.Range("N2").Value = 900

# Row 3: Gerard Moreno -> Alba Redondo
$ws.Range("A3").Value = "Alba Redondo"
$ws.Range("B3").Value = 82
$ws.Range("H3").Value = "Liga F"
$ws.Range("I3").Value = "Levante UD"
$ws.Range("N3").Value = 750

# Row 5: Mewis -> Jensen
$ws.Range("A5").Value = "Jensen"
$ws.Range("B5").Value = 81
$ws.Range("E5").Value = "ST"
$ws.Range("G5").Value = "Norway"
$ws.Range("H5").Value = "Liga F"
$ws.Range("I5").Value = "Real Sociedad"
$ws.Range("O5").Value = 1

# Row 6: Coffey -> Savanier
$ws.Range("A6").Value = "Savanier"
$ws.Range("B6").Value = 80
$ws.Range("E6").Value = "CAM"
$ws.Range("G6").Value = "France"
$ws.Range("H6").Value = "Ligue 1 Uber Eats"
$ws.Range("I6").Value = "Montpellier"

# Row 7: Sørloth -> Mandanda
$ws.Range("A7").Value = "Mandanda"
$ws.Range("E7").Value = "GK"
$ws.Range("G7").Value = "France"
$ws.Range("H7").Value = "Ligue 1 Uber Eats"
$ws.Range("I7").Value = "Stade Rennais FC"
$ws.Range("N7").Value = 550

# Row 8: Webster -> Reguilón
$ws.Range("A8").Value = "Reguilón"
$ws.Range("E8").Value = "LB"
$ws.Range("G8").Value = "Spain"
$ws.Range("I8").Value = "Manchester Utd"
$ws.Range("N8").Value = 450
$ws.Range("O8").Value = 3

# Row 9: Luis Milla -> Chalobah
$ws.Range("A9").Value = "Chalobah"
$ws.Range("B9").Value = 77
$ws.Range("E9").Value = "RB"
$ws.Range("G9").Value = "England"
$ws.Range("H9").Value = "Premier League"
$ws.Range("I9").Value = "Chelsea"
$ws.Range("N9").Value = 450
$ws.Range("O9").Value = 2

# Row 10: Chalobah -> Ferri
$ws.Range("A10").Value = "Ferri"
$ws.Range("E10").Value = "CDM"
$ws.Range("G10").Value = "France"
$ws.Range("H10").Value = "Ligue 1 Uber Eats"
$ws.Range("I10").Value = "Montpellier"
$ws.Range("N10").Value = 500
$ws.Range("O10").Value = 3

# Row 11: Diego Rico -> Anna Torrodà
$ws.Range("A11").Value = "Anna Torrodà"
$ws.Range("B11").Value = 76
$ws.Range("E11").Value = "CDM"
$ws.Range("H11").Value = "Liga F"
$ws.Range("I11").Value = "Levante UD"
